$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I12").Value = 4
